$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "product"
$ws2.Range("A1").Value = "searchkey"
$ws2.Range("B1").Value = "productname"
$ws2.Range("C1").Value = "expectedproductcount"
$ws2.Range("A2").Value = "macbook"
$ws2.Range("B2").Value = "MacBook Pro"
$ws2.Range("C2").Value = 4
$ws2.Range("A3").Value = "macbook"
$ws2.Range("B3").Value = "MacBook Air"
$ws2.Range("C3").Value = 4
$ws2.Range("A4").Value = "imac"
$ws2.Range("B4").Value = "iMac"
$ws2.Range("C4").Value = 3
$ws2.Columns.Item(1).AutoFit() | Out-Null
$ws2.Columns.Item(2).AutoFit() | Out-Null
$ws2.Columns.Item(3).AutoFit() | Out-Null
